$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row labels
$ws.Range("A1").Value = "mx_state"
$ws.Range("B1").Value = "mx_municipality"
$ws.Range("C1").Value = "n_matriculas"
$ws.Range("D1").Value = "pct_matriculas"

# Fix capitalization in municipality names
$ws.Range("B9").Value = "Santa Cruz De Juventino Rosas"
$ws.Range("B16").Value = "Tulancingo De Bravo"

# Remove trailing footer rows (33-37), shifting nothing else since row 32 is blank
$ws.Range("A33:A37").EntireRow.Delete()
